$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 925.46155
$ws.Range("I2").Value = 1269.8889
$ws.Range("J2").Value = 150.5
$ws.Range("K2").Value = 1269.8889
$ws.Range("L2").Value = 150.5
$ws.Range("M2").Value = -1156.8889
$ws.Range("N2").Value = -376.5
$ws.Range("H8").Value = 1138.4286
$ws.Range("I8").Value = 157
$ws.Range("J8").Value = 1531
$ws.Range("K8").Value = 471
$ws.Range("L8").Value = 4593
$ws.Range("M8").Value = -332
$ws.Range("N8").Value = -4871
$ws.Range("H17").Value = 1287.8823
$ws.Range("J17").Value = 1287.8823
$ws.Range("L17").Value = 3863.6469
$ws.Range("N17").Value = -4199.6469
$ws.Range("H38").Value = 1488.762
$ws.Range("I38").Value = 355.41666
$ws.Range("J38").Value = 2999.889
$ws.Range("K38").Value = 1066.24998
$ws.Range("L38").Value = 8999.667000000001
$ws.Range("M38").Value = -694.2499800000001
$ws.Range("N38").Value = -9743.667000000001
$ws.Range("H39").Value = 790.75
$ws.Range("I39").Value = 69.625
$ws.Range("K39").Value = 208.875
$ws.Range("M39").Value = 87.125
$ws.Range("H42").Value = 495.2
$ws.Range("I42").Value = 85.75
$ws.Range("J42").Value = 768.1667
$ws.Range("K42").Value = 257.25
$ws.Range("L42").Value = 2304.5001
$ws.Range("M42").Value = -27.25
$ws.Range("N42").Value = -2764.5001
$ws.Range("H43").Value = 5666.6665
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 5666.6665
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 5666.6665
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -5804.6665
$ws.Range("H58").Value = 19134.695
$ws.Range("J58").Value = 22077.195
$ws.Range("L58").Value = 66231.58499999999
$ws.Range("N58").Value = -66531.58499999999
$ws.Range("H107").Value = 931.3
$ws.Range("I107").Value = 945.4091
$ws.Range("J107").Value = 892.5
$ws.Range("K107").Value = 945.4091
$ws.Range("L107").Value = 892.5
$ws.Range("M107").Value = 974.5909
$ws.Range("N107").Value = -4732.5
$ws.Range("H132").Value = 2704558
$ws.Range("I132").Value = 2986719.5
$ws.Range("J132").Value = 3871.1428
$ws.Range("K132").Value = 8960158.5
$ws.Range("L132").Value = 11613.4284
$ws.Range("M132").Value = -8957628.5
$ws.Range("N132").Value = -16673.4284
$ws.Range("H141").Value = 293877.9
$ws.Range("I141").Value = 1235.7812
$ws.Range("J141").Value = 1854636
$ws.Range("K141").Value = 3707.3436
$ws.Range("L141").Value = 5563908
$ws.Range("M141").Value = 1472.6564
$ws.Range("N141").Value = -5574268

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 35717756
$ws.Range("I2").Value = 41669050
$ws.Range("J2").Value = 10000
$ws.Range("K2").Value = 41669050
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = -41668937
$ws.Range("N2").Value = -10226
$ws.Range("H6").Value = 3882462.2
$ws.Range("I6").Value = 6271500
$ws.Range("J6").Value = 60001.8
$ws.Range("K6").Value = 6271500
$ws.Range("L6").Value = 60001.8
$ws.Range("M6").Value = -6271327
$ws.Range("N6").Value = -60347.8
$ws.Range("H11").Value = 71740.25
$ws.Range("I11").Value = 73476.5
$ws.Range("J11").Value = 70004
$ws.Range("K11").Value = 73476.5
$ws.Range("L11").Value = 70004
$ws.Range("M11").Value = -73332.5
$ws.Range("N11").Value = -70292
$ws.Range("H61").Value = 1361.6724
$ws.Range("I61").Value = 553.6875
$ws.Range("J61").Value = 5240
$ws.Range("K61").Value = 553.6875
$ws.Range("L61").Value = 5240
$ws.Range("M61").Value = -341.6875
$ws.Range("N61").Value = -5664
$ws.Range("H110").Value = 1385.75
$ws.Range("I110").Value = 679.0952
$ws.Range("K110").Value = 679.0952
$ws.Range("M110").Value = 1365.9048
$ws.Range("H116").Value = 35717756
$ws.Range("I116").Value = 41669050
$ws.Range("J116").Value = 10000
$ws.Range("K116").Value = 41669050
$ws.Range("L116").Value = 10000
$ws.Range("M116").Value = -41666756
$ws.Range("N116").Value = -14588
$ws.Range("H136").Value = 1361.6724
$ws.Range("I136").Value = 553.6875
$ws.Range("J136").Value = 5240
$ws.Range("K136").Value = 1661.0625
$ws.Range("L136").Value = 15720
$ws.Range("M136").Value = 888.9375
$ws.Range("N136").Value = -20820

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 35717756
$ws.Range("I3").Value = 41669050
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 41669050
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = -41668936
$ws.Range("N3").Value = -10228
$ws.Range("H12").Value = 270
$ws.Range("I12").Value = 270
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 270
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -102
$ws.Range("N12").ClearContents()
$ws.Range("H94").Value = 874.4194
$ws.Range("I94").Value = 821.85
$ws.Range("J94").Value = 970
$ws.Range("K94").Value = 821.85
$ws.Range("L94").Value = 970
$ws.Range("M94").Value = -370.85
$ws.Range("N94").Value = -1872
$ws.Range("H107").Value = 5433.3335
$ws.Range("I107").Value = 4475
$ws.Range("J107").Value = 6200
$ws.Range("K107").Value = 4475
$ws.Range("L107").Value = 6200
$ws.Range("M107").Value = -2555
$ws.Range("N107").Value = -10040

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 15990.546
$ws.Range("J95").Value = 15990.546
$ws.Range("L95").Value = 15990.546
$ws.Range("N95").Value = -21482.546
$ws.Range("H99").Value = 6337.3335
$ws.Range("I99").Value = 2012
$ws.Range("J99").Value = 8500
$ws.Range("K99").Value = 2012
$ws.Range("L99").Value = 8500
$ws.Range("M99").Value = -514
$ws.Range("N99").Value = -11496
$ws.Range("H105").Value = 2682.9412
$ws.Range("I105").Value = 2678.889
$ws.Range("K105").Value = 2678.889
$ws.Range("M105").Value = -931.8890000000001
$ws.Range("H126").Value = 6337.3335
$ws.Range("I126").Value = 2012
$ws.Range("J126").Value = 8500
$ws.Range("K126").Value = 6036
$ws.Range("L126").Value = 25500
$ws.Range("M126").Value = -3566
$ws.Range("N126").Value = -30440

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 2284.5
$ws.Range("I16").Value = 216.66667
$ws.Range("J16").Value = 2761.6924
$ws.Range("K16").Value = 650.00001
$ws.Range("L16").Value = 8285.0772
$ws.Range("M16").Value = -477.00001
$ws.Range("N16").Value = -8631.0772
$ws.Range("H92").Value = 2783.5715
$ws.Range("I92").Value = 182
$ws.Range("J92").Value = 3217.1667
$ws.Range("K92").Value = 546
$ws.Range("L92").Value = 9651.500100000001
$ws.Range("M92").Value = 702
$ws.Range("N92").Value = -12147.5001
$ws.Range("H131").Value = 1890.1923
$ws.Range("I131").Value = 2101.5386
$ws.Range("J131").Value = 1678.8462
$ws.Range("K131").Value = 6304.6158
$ws.Range("L131").Value = 5036.5386
$ws.Range("M131").Value = -1264.6158
$ws.Range("N131").Value = -15116.5386

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 17613.4
$ws.Range("I3").Value = 1051.5
$ws.Range("J3").Value = 28654.666
$ws.Range("K3").Value = 1051.5
$ws.Range("L3").Value = 28654.666
$ws.Range("M3").Value = -935.5
$ws.Range("N3").Value = -28886.666
$ws.Range("H13").Value = 20235.334
$ws.Range("I13").Value = 412.5
$ws.Range("J13").Value = 27443.637
$ws.Range("K13").Value = 412.5
$ws.Range("L13").Value = 27443.637
$ws.Range("M13").Value = -273.5
$ws.Range("N13").Value = -27721.637
$ws.Range("H113").Value = 7000
$ws.Range("I113").Value = 4000
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 4000
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -1830
$ws.Range("N113").Value = -14340
$ws.Range("H132").Value = 2313
$ws.Range("I132").Value = 1906.1
$ws.Range("K132").Value = 5718.299999999999
$ws.Range("M132").Value = -3188.299999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1363.5217
$ws.Range("I7").Value = 711.35297
$ws.Range("J7").Value = 3211.3333
$ws.Range("K7").Value = 711.35297
$ws.Range("L7").Value = 3211.3333
$ws.Range("M7").Value = -599.35297
$ws.Range("N7").Value = -3435.3333
$ws.Range("H46").Value = 1483
$ws.Range("I46").Value = 480
$ws.Range("J46").Value = 1817.3334
$ws.Range("K46").Value = 480
$ws.Range("L46").Value = 1817.3334
$ws.Range("M46").Value = -292
$ws.Range("N46").Value = -2193.3334
$ws.Range("H126").Value = 1363.5217
$ws.Range("I126").Value = 711.35297
$ws.Range("J126").Value = 3211.3333
$ws.Range("K126").Value = 2134.05891
$ws.Range("L126").Value = 9633.999899999999
$ws.Range("M126").Value = 335.9410899999998
$ws.Range("N126").Value = -14573.9999
$ws.Range("H135").Value = 30306.428
$ws.Range("J135").Value = 30306.428
$ws.Range("L135").Value = 30306.428
$ws.Range("N135").Value = -40446.428

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 120105
$ws.Range("J70").Value = 120105
$ws.Range("L70").Value = 120105
$ws.Range("N70").Value = -120735
$ws.Range("H73").Value = 120105
$ws.Range("J73").Value = 120105
$ws.Range("L73").Value = 120105
$ws.Range("N73").Value = -122289
$ws.Range("H100").Value = 1192.381
$ws.Range("I100").Value = 1224.75
$ws.Range("K100").Value = 2449.5
$ws.Range("M100").Value = -1908.5
$ws.Range("H136").Value = 1045.6863
$ws.Range("I136").Value = 730.7143
$ws.Range("J136").Value = 1734.6875
$ws.Range("K136").Value = 2192.1429
$ws.Range("L136").Value = 5204.0625
$ws.Range("M136").Value = 357.8571000000002
$ws.Range("N136").Value = -10304.0625
